# Update gh-pages to output generated at 456a3b4
# Applies the numeric "want-to-go" counter bumps across the four sheets and
# inserts the new "LoveLive" row into the 演出 (Performance) and 全部类型
# (All types) sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览 (Exhibition) - only F-column (interest-count) bumps
# ---------------------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value  = 872
$wsExhibit.Range("F3").Value  = 13867
$wsExhibit.Range("F4").Value  = 13655
$wsExhibit.Range("F5").Value  = 1056
$wsExhibit.Range("F8").Value  = 607
$wsExhibit.Range("F11").Value = 61
$wsExhibit.Range("F12").Value = 773
$wsExhibit.Range("F13").Value = 2155
$wsExhibit.Range("F15").Value = 95
$wsExhibit.Range("F17").Value = 131
$wsExhibit.Range("F19").Value = 540
$wsExhibit.Range("F20").Value = 440
$wsExhibit.Range("F21").Value = 420
$wsExhibit.Range("F24").Value = 846
$wsExhibit.Range("F25").Value = 101

# ---------------------------------------------------------------------
# Sheet 2: 演出 (Performance) - F-column bumps + new row inserted at 14
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value  = 23
$wsShow.Range("F7").Value  = 1559
$wsShow.Range("F12").Value = 72

# Insert a brand-new row before the current row 14 (孟京辉 ...), pushing it
# (and everything below) down by one.
$wsShow.Rows.Item(14).Insert()

# Restore the index-column formatting (bold/centered/bordered) on the new
# row 14 that the insert otherwise leaves unformatted.
$wsShow.Range("A13").Copy()
$wsShow.Range("A14").PasteSpecial(-4122)

$wsShow.Range("A14").Value = 13
$wsShow.Range("B14").NumberFormat = "@"
$wsShow.Range("B14").Value = "2024-08-24"
$wsShow.Range("C14").Value = "广州·LoveLive！电视动画播放十周年纪念巡演"
$wsShow.Range("D14").Value = "机场路1733号 久米空间LIVEHOUSE"
$wsShow.Range("E14").Value = "2024.08.24 12:30-08.25 18:30"
$wsShow.Range("F14").Value = 2
$wsShow.Range("G14").Value = 580
$wsShow.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=86959"
$wsShow.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202406/apzqBc5d1717661406596.jpeg"

# Renumber the (now pushed-down) old row 14 -> row 15's index column.
$wsShow.Range("A15").Value = 14

# ---------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life) - only F-column bumps
# ---------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 223
$wsLocal.Range("F3").Value = 118

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (All types) - F-column bumps + new row inserted at 42
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 223
$wsAll.Range("F3").Value  = 872
$wsAll.Range("F4").Value  = 13867
$wsAll.Range("F5").Value  = 13655
$wsAll.Range("F6").Value  = 1056
$wsAll.Range("F9").Value  = 607
$wsAll.Range("F12").Value = 61
$wsAll.Range("F13").Value = 773
$wsAll.Range("F15").Value = 23
$wsAll.Range("F16").Value = 2155
$wsAll.Range("F18").Value = 95
$wsAll.Range("F20").Value = 131
$wsAll.Range("F24").Value = 118
$wsAll.Range("F25").Value = 118
$wsAll.Range("F26").Value = 540
$wsAll.Range("F27").Value = 440
$wsAll.Range("F28").Value = 420
$wsAll.Range("F31").Value = 846
$wsAll.Range("F33").Value = 1559
$wsAll.Range("F38").Value = 101
$wsAll.Range("F39").Value = 72

# Insert a brand-new row before the current row 42 (孟京辉 ...), pushing it
# (and everything below) down by one.
$wsAll.Rows.Item(42).Insert()

# Restore the index-column formatting (bold/centered/bordered) on the new
# row 42 that the insert otherwise leaves unformatted.
$wsAll.Range("A41").Copy()
$wsAll.Range("A42").PasteSpecial(-4122)

$wsAll.Range("A42").Value = 41
$wsAll.Range("B42").NumberFormat = "@"
$wsAll.Range("B42").Value = "2024-08-24"
$wsAll.Range("C42").Value = "广州·LoveLive！电视动画播放十周年纪念巡演"
$wsAll.Range("D42").Value = "机场路1733号 久米空间LIVEHOUSE"
$wsAll.Range("E42").Value = "2024.08.24 12:30-08.25 18:30"
$wsAll.Range("F42").Value = 2
$wsAll.Range("G42").Value = 580
$wsAll.Range("H42").Value = "https://show.bilibili.com/platform/detail.html?id=86959"
$wsAll.Range("I42").Value = "//i1.hdslb.com/bfs/openplatform/202406/apzqBc5d1717661406596.jpeg"

# Renumber the (now pushed-down) old row 42 -> row 43's index column.
$wsAll.Range("A43").Value = 42
